# Refresh cryptocurrency price/volume data (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need the Text number format
# applied first so Excel stores them as text (matching the source data feed),
# rather than auto-converting them to floating point numbers.
$textForceCells = @("D5", "D6", "D7", "D9", "D10", "D13", "D17", "D19", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D38", "D41", "D43", "D44", "D47", "D48", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '51.420.16'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '2.774.74'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '354.48'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').Value = '108.00'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').Value = '0.549'
$ws.Range('E7').Value = '  -2.67%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.587'
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').Value = '39.69'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('E11').Value = '  +3.30%  '
$ws.Range('D13').Value = '0.0833'
$ws.Range('E13').Value = '  -2.29%  '
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = '3.209.16'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').Value = '2.774.41'
$ws.Range('E16').Value = '  -2.82%  '
$ws.Range('D17').Value = '0.925'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '51.408.83'
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('D19').Value = '7.64'
$ws.Range('E19').Value = '  +2.96%  '
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('D21').Value = '13.20'
$ws.Range('E21').Value = '  +1.21%  '
$ws.Range('D22').Value = '0.0₃0963'
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('D23').Value = '69.77'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').Value = '265.82'
$ws.Range('E24').Value = '  -3.05%  '
$ws.Range('D25').Value = '2.70'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '26.01'
$ws.Range('E27').Value = '  -2.01%  '
$ws.Range('D29').Value = '10.22'
$ws.Range('E29').Value = '  +0.73%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '36.55'
$ws.Range('E30').Value = '  +7.54%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '2.20'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').Value = '6.16'
$ws.Range('E32').Value = '  +8.12%  '
$ws.Range('D33').Value = '51.79'
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('D34').Value = '0.0444'
$ws.Range('E34').Value = '  -4.33%  '
$ws.Range('D35').Value = '5.53'
$ws.Range('E35').Value = '  +5.50%  '
$ws.Range('E36').Value = '  -2.32%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').Value = '18.25'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('E39').Value = '  -2.58%  '
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('D41').Value = '2.52'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('D43').Value = '120.87'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').Value = '22.04'
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('D46').Value = '2.099.57'
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '2.32'
$ws.Range('E47').Value = '  +4.01%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '3.25'
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('E49').Value = '  -4.69%  '
$ws.Range('D50').Value = '0.905'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('E51').Value = '  +7.71%  '
